$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# B1: "Dollar" -> "Total Lobbies"
$ws.Range("B1").Value = "Total Lobbies"

# New headers for the two added columns
$ws.Range("C1").Value = "Lobbies to Democrats"
$ws.Range("D1").Value = "Lobbies to Republicans "

# Match the header formatting used by A1/B1 (bold font), but without the
# wrap/center alignment that A1/B1 use.
$ws.Range("C1:D1").Font.Bold = $true

# --- Add the "Lobbies to Democrats" / "Lobbies to Republicans" data ---
# Data only exists for even years (every other row).
$ws.Range("C2").Value = 707800
$ws.Range("D2").Value = 1300000

$ws.Range("C4").Value = 809200
$ws.Range("D4").Value = 2100000

$ws.Range("C6").Value = 918200
$ws.Range("D6").Value = 2200000

$ws.Range("C8").Value = 1000000
$ws.Range("D8").Value = 2400000

$ws.Range("C10").Value = 1100000
$ws.Range("D10").Value = 2900000

$ws.Range("C12").Value = 2000000
$ws.Range("D12").Value = 3100000

$ws.Range("C14").Value = 1600000
$ws.Range("D14").Value = 2400000

$ws.Range("C16").Value = 999000
$ws.Range("D16").Value = 3700000

$ws.Range("C18").Value = 791500
$ws.Range("D18").Value = 2800000

$ws.Range("C20").Value = 943400
$ws.Range("D20").Value = 3000000

$ws.Range("C22").Value = 821200
$ws.Range("D22").Value = 2400000
